$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(2, 0, 2, 1, 1, 1, 0, 1, 0, 2, 0, 2, 2, 0, 1, 1, 1)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $values[$i]
}
